# Natmi following Dr Hou advice
# Rebuild the LR-pair table for Nodal-Acvr2b to include "ECs" as a sending
# cluster (previously only FAPs/sCs were sending clusters), recomputed with
# updated specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Nodal", "Acvr2b", "ECs", 2, 0.6666666666666666, 0.3514346666666666, 1.054304, 0.164942595720082, 0.1649425957200821, 3, 1, 0.5292956666666667, 1.587887, 0.135651968140022, 0.1356519681400219, 0.1860128461831111, 1.674115615648, 0.02237478773955309, 0.02237478773955309),
    @("ECs", "Nodal", "Acvr2b", "FAPs", 2, 0.6666666666666666, 0.3514346666666666, 1.054304, 0.164942595720082, 0.1649425957200821, 3, 1, 1.362890666666667, 4.088672000000001, 0.3492921120199358, 0.3492921120199358, 0.4789670271431111, 4.310703244288001, 0.05761314762111788, 0.0576131476211179),
    @("ECs", "Nodal", "Acvr2b", "sCs", 2, 0.6666666666666666, 0.3514346666666666, 1.054304, 0.164942595720082, 0.1649425957200821, 3, 1, 2.009678666666666, 6.029036, 0.5150559198400423, 0.5150559198400423, 0.7062707523271109, 6.356436770943999, 0.08495466035941109, 0.0849546603594111),
    @("FAPs", "Nodal", "Acvr2b", "ECs", 3, 1, 1.395965, 4.187895, 0.6551832032346963, 0.6551832032346963, 3, 1, 0.5292956666666667, 1.587887, 0.135651968140022, 0.1356519681400219, 0.7388782253183335, 6.649904027865, 0.08887689101107055, 0.08887689101107052),
    @("FAPs", "Nodal", "Acvr2b", "FAPs", 3, 1, 1.395965, 4.187895, 0.6551832032346963, 0.6551832032346963, 3, 1, 1.362890666666667, 4.088672000000001, 0.3492921120199358, 0.3492921120199358, 1.902547669493334, 17.12292902544, 0.2288503248178339, 0.2288503248178339),
    @("FAPs", "Nodal", "Acvr2b", "sCs", 3, 1, 1.395965, 4.187895, 0.6551832032346963, 0.6551832032346963, 3, 1, 2.009678666666666, 6.029036, 0.5150559198400423, 0.5150559198400423, 2.805441079913333, 25.24896971922, 0.3374559874057919, 0.3374559874057919),
    @("sCs", "Nodal", "Acvr2b", "ECs", 3, 1, 0.3832486666666666, 1.149746, 0.1798742010452218, 0.1798742010452218, 3, 1, 0.5292956666666667, 1.587887, 0.135651968140022, 0.1356519681400219, 0.2028518585224444, 1.825666726702, 0.02440028938939832, 0.02440028938939832),
    @("sCs", "Nodal", "Acvr2b", "FAPs", 3, 1, 0.3832486666666666, 1.149746, 0.1798742010452218, 0.1798742010452218, 3, 1, 1.362890666666667, 4.088672000000001, 0.3492921120199358, 0.3492921120199358, 0.5223260308124444, 4.700934277312, 0.06282863958098406, 0.06282863958098406),
    @("sCs", "Nodal", "Acvr2b", "sCs", 3, 1, 0.3832486666666666, 1.149746, 0.1798742010452218, 0.1798742010452218, 3, 1, 2.009678666666666, 6.029036, 0.5150559198400423, 0.5150559198400423, 0.7702066694284443, 6.931860024855999, 0.09264527207483939, 0.09264527207483939)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

